$p = $ppt.ActivePresentation

$oldDate = "8/24/23"
$newDate = "2/19/24"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Characters(1, $tr.Length).Text = $newDate
            }
        }
    }
}

# Update the date placeholder on the slide master.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Update the date placeholder on every slide layout.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Touch up the experimenter instructions on slide 1.
$slide = $p.Slides.Item(1)
$shape = $slide.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange
$oldPhrase = "Two more sounds will be played. "
$newPhrase = "Another two sounds will be played. "
$para1 = $tr.Paragraphs(1, 1)
$para1Text = $para1.Text.TrimEnd("`r", "`n")
if ($para1Text -eq $oldPhrase) {
    $tr.Characters($para1.Start, $oldPhrase.Length).Text = $newPhrase
}
